$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46081 -> 46082) for every data row (rows 2-95).
for ($row = 2; $row -le 95; $row++) {
    $ws.Cells.Item($row, 3).Value = 46082
}
